# Generate Report for Handback
# Marks the handback as complete (status text), records the handback
# timestamp, and fills in the "Latest Target File" / "Latest Handback File"
# / "Latest Handback DateTime" columns (with hyperlinks on the target file
# column) for both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Status text for both rows/languages moves from "Ready for handoff" to
# "Handed back: in sync with en-US" - this text is shared between the
# Overview roll-up columns (E/F) and each language sheet's "Status"
# column (C), so update every occurrence consistently.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"
$overview.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$overview.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet: fill Latest Target File / Latest Handback File / DateTime
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn.Range("I2").Value = "52c54f0c-3e0c-4864-94c2-7118823c50dd.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f5a2e6b57cc1c9a94779bd70b1d629ea60acdc3/e2e/52c54f0c-3e0c-4864-94c2-7118823c50dd.md", "", "", "52c54f0c-3e0c-4864-94c2-7118823c50dd.md") | Out-Null
$zhcn.Range("J2").Value = "52c54f0c-3e0c-4864-94c2-7118823c50dd.6086963e0aadb8498a285965140a6b5e0fb14e3b.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-06 03:33:58"

$zhcn.Range("I3").Value = "b7e88faa-6dbd-4f8c-8fde-a41bee98a83d.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f5a2e6b57cc1c9a94779bd70b1d629ea60acdc3/e2e/b7e88faa-6dbd-4f8c-8fde-a41bee98a83d.md", "", "", "b7e88faa-6dbd-4f8c-8fde-a41bee98a83d.md") | Out-Null
$zhcn.Range("J3").Value = "b7e88faa-6dbd-4f8c-8fde-a41bee98a83d.feafb79ab31951430c8076a4e1deeb37e0ed07e2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-06 03:33:58"

$zhcn.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$zhcn.Columns.Item(9).EntireColumn.AutoFit() | Out-Null
$zhcn.Columns.Item(10).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: fill Latest Target File / Latest Handback File / DateTime
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Range("I2").Value = "52c54f0c-3e0c-4864-94c2-7118823c50dd.md"
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f5a2e6b57cc1c9a94779bd70b1d629ea60acdc3/e2e/52c54f0c-3e0c-4864-94c2-7118823c50dd.md", "", "", "52c54f0c-3e0c-4864-94c2-7118823c50dd.md") | Out-Null
$dede.Range("J2").Value = "52c54f0c-3e0c-4864-94c2-7118823c50dd.6086963e0aadb8498a285965140a6b5e0fb14e3b.de-de.xlf"
$dede.Range("K2").Value = "2016-09-06 03:34:20"

$dede.Range("I3").Value = "b7e88faa-6dbd-4f8c-8fde-a41bee98a83d.md"
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f5a2e6b57cc1c9a94779bd70b1d629ea60acdc3/e2e/b7e88faa-6dbd-4f8c-8fde-a41bee98a83d.md", "", "", "b7e88faa-6dbd-4f8c-8fde-a41bee98a83d.md") | Out-Null
$dede.Range("J3").Value = "b7e88faa-6dbd-4f8c-8fde-a41bee98a83d.feafb79ab31951430c8076a4e1deeb37e0ed07e2.de-de.xlf"
$dede.Range("K3").Value = "2016-09-06 03:34:20"

$dede.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$dede.Columns.Item(9).EntireColumn.AutoFit() | Out-Null
$dede.Columns.Item(10).EntireColumn.AutoFit() | Out-Null
